# BuildHeap Youtube presentation started ArrayHeap fix for pull mehtod when
# pulling last element. Optimize minHeapifyRoot and maxHeapifyRoot. Change
# the way a swap is done in maxHeapifyRoot.
#
# The recorded OOXML diff for this commit only touches the cached
# "datetimeFigureOut" date fields that live on the slide master, every
# slide layout, and the notes master (the deck's slides themselves don't
# carry a date placeholder) -- each one flips from 3/2/2016 to 3/3/2016.
# Update every "Date Placeholder" shape we can find across those
# containers.

$p = $ppt.ActivePresentation

$oldDate = "3/2/2016"
$newDate = "3/3/2016"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1) Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2) Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lyt = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $lyt.Shapes
}

# 3) Notes master.
$nm = $p.NotesMaster
Update-DatePlaceholder $nm.Shapes
